# Fruta / hortaliza, semanal
# Insert a new daily price record as row 27 in the "Macroferia Regional de
# Talca - Arándano (blue)" sheet. All existing rows from 27 downward shift
# down by one (Excel's normal Insert-row behaviour), and the new row is
# populated with the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 27..99 down to 28..100, opening up a blank row 27.
$ws.Rows(27).Insert()

# Fill the newly opened row 27 with the new record.
$ws.Cells.Item(27, 1).Value  = 5
$ws.Cells.Item(27, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(27, 3).Value  = "Maule"
$ws.Cells.Item(27, 4).Value  = 44925
$ws.Cells.Item(27, 5).Value  = 7
$ws.Cells.Item(27, 6).Value  = "Fruta"
$ws.Cells.Item(27, 7).Value  = 100101
$ws.Cells.Item(27, 8).Value  = "Berries"
$ws.Cells.Item(27, 9).Value  = 100101001
$ws.Cells.Item(27, 10).Value = "Arándano (blue)"
$ws.Cells.Item(27, 11).Value = "Sin especificar"
$ws.Cells.Item(27, 12).Value = "Primera"
$ws.Cells.Item(27, 13).Value = 150
$ws.Cells.Item(27, 14).Value = 3000
$ws.Cells.Item(27, 15).Value = 3000
$ws.Cells.Item(27, 16).Value = 3000
$ws.Cells.Item(27, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(27, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(27, 19).Value = 1500
$ws.Cells.Item(27, 20).Value = 2
